$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; this shifts the existing rows 7..50 down to 8..51
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new weekly data point.
# (Same Mercado/Region/Categoria/Variedad/Calidad/Unidad/Origen/Kg metadata
# as the surrounding rows; new Fecha + Volumen + Precios.)
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 45083
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 100112043
$ws.Range("G7").Value = "Pepino dulce"
$ws.Range("H7").Value = "Cultivar IV Región"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = 14000
$ws.Range("N7").Value = "$/bandeja 18 kilos"
$ws.Range("O7").Value = "Provincia de Limarí"
$ws.Range("P7").Value = 778
$ws.Range("Q7").Value = 18
$ws.Range("R7").Value = "Hortaliza"
